$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 0.7536891456382829
$ws.Range("D2").Value = 0.1689636070299383
$ws.Range("E2").Value = "norm_qa_fix_dispersion_mean"
$ws.Range("F2").Value = "'2.08e-09"
$ws.Range("F2").ClearFormats()

# Row 3
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 0.4281139994473667
$ws.Range("D3").Value = 0.13074018657193
$ws.Range("E3").Value = "norm_qa_fix_dispersion_mean"
$ws.Range("F3").Value = "'7.64e-08"
$ws.Range("F3").ClearFormats()

# Row 4
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 0.5419255782647949
$ws.Range("D4").Value = 0.1435993001297624
$ws.Range("E4").Value = "norm_coldread_coverage_line_%"
$ws.Range("F4").Value = "'2.87e-05"
$ws.Range("F4").ClearFormats()
